$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 41, shifting existing rows 41:135 down to 43:137
$ws.Rows("41:42").Insert()

# --- Populate new row 41 ---
$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C41").Value = "Los Lagos"
$ws.Range("D41").Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100104
$ws.Range("H41").Value = "Frutos de pepita"
$ws.Range("I41").Value = 100104005
$ws.Range("J41").Value = "Pera"
$ws.Range("K41").Value = "Packham's Triumph"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 200
$ws.Range("N41").Value = 16000
$ws.Range("O41").Value = 16000
$ws.Range("P41").Value = 16000
$ws.Range("Q41").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R41").Value = "Región de O'Higgins"
$ws.Range("S41").Value = 1067
$ws.Range("T41").Value = 15

# --- Populate new row 42 ---
$ws.Range("A42").Value = 4
$ws.Range("B42").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C42").Value = "Los Lagos"
$ws.Range("D42").Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E42").Value = 10
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100104
$ws.Range("H42").Value = "Frutos de pepita"
$ws.Range("I42").Value = 100104005
$ws.Range("J42").Value = "Pera"
$ws.Range("K42").Value = "Packham's Triumph"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 13000
$ws.Range("O42").Value = 13000
$ws.Range("P42").Value = 13000
$ws.Range("Q42").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R42").Value = "Región de O'Higgins"
$ws.Range("S42").Value = 867
$ws.Range("T42").Value = 15
